$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header F1: shared-string text "C" -> "B", and center-aligned like the data below it
$ws.Range("F1").Value = "B"
$ws.Range("F1").HorizontalAlignment = -4108

# F2:F8: replace each "ohms" figure with its reciprocal ("per-unit" style figure)
foreach ($r in 2..8) {
    $cell = $ws.Cells.Item($r, 6)
    $old = $cell.Value2
    $cell.Value = 1 / $old
}

# New empty, centre-column-formatted cells I3:I8 (same style as G3:G8) extend the used range to column I
$ws.Range("G3:G8").Copy()
$ws.Range("I3:I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the visible selection to match the author's last cursor position
$ws.Range("K10").Select()
